$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 217
$ws.Range("F3").Value = 1386
$ws.Range("F4").Value = 19515
$ws.Range("F6").Value = 301
$ws.Range("F8").Value = 6
$ws.Range("F9").Value = 7385
$ws.Range("F10").Value = 482
$ws.Range("F11").Value = 719
$ws.Range("F13").Value = 33
$ws.Range("F14").Value = 148
$ws.Range("F15").Value = 101
$ws.Range("F17").Value = 181
$ws.Range("F19").Value = 362
$ws.Range("F22").Value = 44
$ws.Range("F23").Value = 48
$ws.Range("F24").Value = 57
$ws.Range("F25").Value = 307
$ws.Range("F26").Value = 1066
$ws.Range("F28").Value = 7
$ws.Range("F29").Value = 163
$ws.Range("F32").Value = 46
$ws.Range("F33").Value = 964
$ws.Range("F35").Value = 84
$ws.Range("F36").Value = 10
$ws.Range("F37").Value = 12481
$ws.Range("F38").Value = 1317
$ws.Range("F39").Value = 56
$ws.Range("F40").Value = 10
$ws.Range("F43").Value = 334
$ws.Range("F44").Value = 0
$ws.Range("F45").Value = 318

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 2

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 217
$ws.Range("F3").Value = 1386
$ws.Range("F4").Value = 19515
$ws.Range("F5").Value = 781
$ws.Range("F6").Value = 301
$ws.Range("F9").Value = 7385
$ws.Range("F10").Value = 482
$ws.Range("F11").Value = 719
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = 33
$ws.Range("F14").Value = 148
$ws.Range("F15").Value = 101
$ws.Range("F17").Value = 181
$ws.Range("F18").Value = 1327
$ws.Range("F19").Value = 362
$ws.Range("F20").Value = 66
$ws.Range("F21").Value = 675
$ws.Range("F22").Value = 44
$ws.Range("F23").Value = 48
$ws.Range("F24").Value = 57
$ws.Range("F25").Value = 307
$ws.Range("F26").Value = 1066
$ws.Range("F28").Value = 7
$ws.Range("F31").Value = 551
$ws.Range("F32").Value = 2
$ws.Range("F37").Value = 84
$ws.Range("F38").Value = 10
$ws.Range("F39").Value = 12481
$ws.Range("F40").Value = 1317
$ws.Range("F41").Value = 56
$ws.Range("F43").Value = 52
$ws.Range("F45").Value = 334
$ws.Range("F46").Value = 0
